$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 23: PFM / Transaction Details / SCPAccounts / GET / endpoint / issue description
$ws.Range("A23").Value = "PFM"
$ws.Range("B23").Value = "Transaction Details"
$ws.Range("C23").Value = "SCPAccounts"
$ws.Range("D23").Value = "GET"
$ws.Range("E23").Value = "/scp/account/transactions "
$ws.Range("F23").Value = "ChildTransactiosns are what   basis for the API:- account/transactions in SCPAccounts.yaml"

# New rows 24 and 25 just carry the PFM stream label forward
$ws.Range("A24").Value = "PFM"
$ws.Range("A25").Value = "PFM"

# Update the active selection to match the new view state
$ws.Range("C25").Select()
